$d = $word.ActiveDocument
$p62 = $d.Paragraphs.Item(62)
$rEnd = $p62.Range
$rEnd.InsertParagraphAfter()
$p63 = $d.Paragraphs.Item(63)
$pr = $p63.Range
$startPos = $pr.Start
$fullText = "plot" + "(st.rel[st.rel" + "`$" + "period" + "==" + "`"[1993,2003)`"" + ",])" + [char]11 + "lines" + "(st.rel[st.rel" + "`$" + "period" + "==" + "`"[2003,2013)`"" + ",]," + "col=" + "`"red`"" + ")" + [char]11 + "legend" + "(" + "`"topright`"" + "," + "c" + "(" + "`"1993-2002`"" + "," + "`"2003-2012`"" + ")," + "col=" + "c" + "(" + "`"black`"" + "," + "`"red`"" + ")," + "lty=" + "1" + ")"
$pr.Text = $fullText

$pos = $startPos
$seg = $d.Range($pos, $pos + 4)  # "plot"
$seg.Style = "KeywordTok"
$pos = $pos + 4
$seg = $d.Range($pos, $pos + 14)  # "(st.rel[st.rel"
$seg.Style = "NormalTok"
$pos = $pos + 14
$seg = $d.Range($pos, $pos + 1)  # "`$"
$seg.Style = "OperatorTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 6)  # "period"
$seg.Style = "NormalTok"
$pos = $pos + 6
$seg = $d.Range($pos, $pos + 2)  # "=="
$seg.Style = "OperatorTok"
$pos = $pos + 2
$seg = $d.Range($pos, $pos + 13)  # "`"[1993,2003)`""
$seg.Style = "StringTok"
$pos = $pos + 13
$seg = $d.Range($pos, $pos + 3)  # ",])"
$seg.Style = "NormalTok"
$pos = $pos + 3
$pos = $pos + 1  # BR
$seg = $d.Range($pos, $pos + 5)  # "lines"
$seg.Style = "KeywordTok"
$pos = $pos + 5
$seg = $d.Range($pos, $pos + 14)  # "(st.rel[st.rel"
$seg.Style = "NormalTok"
$pos = $pos + 14
$seg = $d.Range($pos, $pos + 1)  # "`$"
$seg.Style = "OperatorTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 6)  # "period"
$seg.Style = "NormalTok"
$pos = $pos + 6
$seg = $d.Range($pos, $pos + 2)  # "=="
$seg.Style = "OperatorTok"
$pos = $pos + 2
$seg = $d.Range($pos, $pos + 13)  # "`"[2003,2013)`""
$seg.Style = "StringTok"
$pos = $pos + 13
$seg = $d.Range($pos, $pos + 3)  # ",],"
$seg.Style = "NormalTok"
$pos = $pos + 3
$seg = $d.Range($pos, $pos + 4)  # "col="
$seg.Style = "DataTypeTok"
$pos = $pos + 4
$seg = $d.Range($pos, $pos + 5)  # "`"red`""
$seg.Style = "StringTok"
$pos = $pos + 5
$seg = $d.Range($pos, $pos + 1)  # ")"
$seg.Style = "NormalTok"
$pos = $pos + 1
$pos = $pos + 1  # BR
$seg = $d.Range($pos, $pos + 6)  # "legend"
$seg.Style = "KeywordTok"
$pos = $pos + 6
$seg = $d.Range($pos, $pos + 1)  # "("
$seg.Style = "NormalTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 10)  # "`"topright`""
$seg.Style = "StringTok"
$pos = $pos + 10
$seg = $d.Range($pos, $pos + 1)  # ","
$seg.Style = "NormalTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 1)  # "c"
$seg.Style = "KeywordTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 1)  # "("
$seg.Style = "NormalTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 11)  # "`"1993-2002`""
$seg.Style = "StringTok"
$pos = $pos + 11
$seg = $d.Range($pos, $pos + 1)  # ","
$seg.Style = "NormalTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 11)  # "`"2003-2012`""
$seg.Style = "StringTok"
$pos = $pos + 11
$seg = $d.Range($pos, $pos + 2)  # "),"
$seg.Style = "NormalTok"
$pos = $pos + 2
$seg = $d.Range($pos, $pos + 4)  # "col="
$seg.Style = "DataTypeTok"
$pos = $pos + 4
$seg = $d.Range($pos, $pos + 1)  # "c"
$seg.Style = "KeywordTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 1)  # "("
$seg.Style = "NormalTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 7)  # "`"black`""
$seg.Style = "StringTok"
$pos = $pos + 7
$seg = $d.Range($pos, $pos + 1)  # ","
$seg.Style = "NormalTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 5)  # "`"red`""
$seg.Style = "StringTok"
$pos = $pos + 5
$seg = $d.Range($pos, $pos + 2)  # "),"
$seg.Style = "NormalTok"
$pos = $pos + 2
$seg = $d.Range($pos, $pos + 4)  # "lty="
$seg.Style = "DataTypeTok"
$pos = $pos + 4
$seg = $d.Range($pos, $pos + 1)  # "1"
$seg.Style = "DecValTok"
$pos = $pos + 1
$seg = $d.Range($pos, $pos + 1)  # ")"
$seg.Style = "NormalTok"
$pos = $pos + 1
